# Update "Pais" sheet with the latest COVID-19 country snapshot
# (countries re-ranked + Spain provincias data refresh, 19 May 2020 11:35)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp in the title row ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 11:35"

# --- Row data updates -------------------------------------------------
# Each entry: row number, country name, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @(
    @(4,   "Estados Unidos", 1550539, 245, 356383, 1102171, 0, 4, 91985),
    @(19,  "Belgica",        55791,   232, 14687,  31996,   0, 28, 9108),
    @(38,  "Israel",         16650,   7,   13299,  3074,    0, 1,  277),
    @(54,  "Afganistan",     7653,    581, 850,    6625,    0, 5,  178),
    @(55,  "Argelia",        7201,    0,   3625,   3021,    0, 0,  555),
    @(56,  "Barein",         7184,    0,   2931,   4241,    0, 0,  12),
    @(196, "Belice",         18,      0,   16,     0,       0, 0,  2),
    @(197, "Nueva Caledonia",18,      0,   18,     0,       0, 0,  0),
    @(209, "Groenlandia",    11,      0,   11,     0,       0, 0,  0),
    @(210, "Montserrat",     11,      0,   10,     0,       0, 0,  1),
    @(211, "Seychelles",     11,      0,   11,     0,       0, 0,  0),
    @(215, "Bonaire, San Eustaquio y Saba", 6, 0, 6, 0, 0, 0, 0),
    @(216, "San Bartolome",  6,       0,   6,      0,       0, 0,  0)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 1).Value = $u[1]
    $ws.Cells.Item($row, 2).Value = $u[2]
    $ws.Cells.Item($row, 3).Value = $u[3]
    $ws.Cells.Item($row, 4).Value = $u[4]
    $ws.Cells.Item($row, 5).Value = $u[5]
    $ws.Cells.Item($row, 6).Value = $u[6]
    $ws.Cells.Item($row, 7).Value = $u[7]
    $ws.Cells.Item($row, 8).Value = $u[8]
}
